# Auto-generated edit script applying scheduled-runner price/profit updates
# to the Tiamat_Profits workbook (per-job sheets ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR).
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(17, 8).Value = 449566.12   # H17
$ws.Cells.Item(17, 10).Value = 449566.12   # J17
$ws.Cells.Item(17, 12).Value = 1348698.36   # L17
$ws.Cells.Item(17, 14).Value = -1349034.36   # N17
$ws.Cells.Item(20, 8).Value = 8157.143   # H20
$ws.Cells.Item(20, 9).Value = 1683.3334   # I20
$ws.Cells.Item(20, 10).Value = 47000   # J20
$ws.Cells.Item(20, 11).Value = 1683.3334   # K20
$ws.Cells.Item(20, 12).Value = 47000   # L20
$ws.Cells.Item(20, 13).Value = -1453.3334   # M20
$ws.Cells.Item(20, 14).Value = -47460   # N20
$ws.Cells.Item(35, 8).Value = 8157.143   # H35
$ws.Cells.Item(35, 9).Value = 1683.3334   # I35
$ws.Cells.Item(35, 10).Value = 47000   # J35
$ws.Cells.Item(35, 11).Value = 1683.3334   # K35
$ws.Cells.Item(35, 12).Value = 47000   # L35
$ws.Cells.Item(35, 13).Value = -1304.3334   # M35
$ws.Cells.Item(35, 14).Value = -47758   # N35
$ws.Cells.Item(116, 8).Value = 5430.1953   # H116
$ws.Cells.Item(116, 9).Value = 6695.2173   # I116
$ws.Cells.Item(116, 10).Value = 3813.7778   # J116
$ws.Cells.Item(116, 11).Value = 6695.2173   # K116
$ws.Cells.Item(116, 12).Value = 3813.7778   # L116
$ws.Cells.Item(116, 13).Value = -3253.2173   # M116
$ws.Cells.Item(116, 14).Value = -10697.7778   # N116
$ws.Cells.Item(137, 8).Value = 4265.3823   # H137
$ws.Cells.Item(137, 9).Value = 1103.25   # I137
$ws.Cells.Item(137, 10).Value = 4687   # J137
$ws.Cells.Item(137, 11).Value = 3309.75   # K137
$ws.Cells.Item(137, 12).Value = 14061   # L137
$ws.Cells.Item(137, 13).Value = -759.75   # M137
$ws.Cells.Item(137, 14).Value = -19161   # N137

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(61, 8).Value = 2456.2727   # H61
$ws.Cells.Item(61, 9).Value = 1651.9   # I61
$ws.Cells.Item(61, 11).Value = 1651.9   # K61
$ws.Cells.Item(61, 13).Value = -1439.9   # M61
$ws.Cells.Item(74, 8).Value = 43827.457   # H74
$ws.Cells.Item(74, 9).Value = 74044.36   # I74
$ws.Cells.Item(74, 10).Value = 1523.8   # J74
$ws.Cells.Item(74, 11).Value = 74044.36   # K74
$ws.Cells.Item(74, 12).Value = 1523.8   # L74
$ws.Cells.Item(74, 13).Value = -73170.36   # M74
$ws.Cells.Item(74, 14).Value = -3271.8   # N74
$ws.Cells.Item(77, 8).Value = 43827.457   # H77
$ws.Cells.Item(77, 9).Value = 74044.36   # I77
$ws.Cells.Item(77, 10).Value = 1523.8   # J77
$ws.Cells.Item(77, 11).Value = 370221.8   # K77
$ws.Cells.Item(77, 12).Value = 7619   # L77
$ws.Cells.Item(77, 13).Value = -365853.8   # M77
$ws.Cells.Item(77, 14).Value = -16355   # N77
$ws.Cells.Item(88, 8).Value = 35666.332   # H88
$ws.Cells.Item(88, 9).Value = 2500   # I88
$ws.Cells.Item(88, 10).Value = 52249.5   # J88
$ws.Cells.Item(88, 11).Value = 2500   # K88
$ws.Cells.Item(88, 12).Value = 52249.5   # L88
$ws.Cells.Item(88, 13).Value = -2094   # M88
$ws.Cells.Item(88, 14).Value = -53061.5   # N88
$ws.Cells.Item(91, 8).Value = 35666.332   # H91
$ws.Cells.Item(91, 9).Value = 2500   # I91
$ws.Cells.Item(91, 10).Value = 52249.5   # J91
$ws.Cells.Item(91, 11).Value = 2500   # K91
$ws.Cells.Item(91, 12).Value = 52249.5   # L91
$ws.Cells.Item(91, 13).Value = -1096   # M91
$ws.Cells.Item(91, 14).Value = -55057.5   # N91
$ws.Cells.Item(110, 8).Value = 636.5   # H110
$ws.Cells.Item(110, 9).Value = 505.7143   # I110
$ws.Cells.Item(110, 10).Value = 819.6   # J110
$ws.Cells.Item(110, 11).Value = 505.7143   # K110
$ws.Cells.Item(110, 12).Value = 819.6   # L110
$ws.Cells.Item(110, 13).Value = 1539.2857   # M110
$ws.Cells.Item(110, 14).Value = -4909.6   # N110
$ws.Cells.Item(122, 8).Value = 1289.5   # H122
$ws.Cells.Item(122, 9).Value = 1031.6666   # I122
$ws.Cells.Item(122, 10).Value = 1400   # J122
$ws.Cells.Item(122, 11).Value = 3094.9998   # K122
$ws.Cells.Item(122, 12).Value = 4200   # L122
$ws.Cells.Item(122, 13).Value = -644.9998000000001   # M122
$ws.Cells.Item(122, 14).Value = -9100   # N122
$ws.Cells.Item(132, 8).Value = 2180079.2   # H132
$ws.Cells.Item(132, 9).Value = 2464781   # I132
$ws.Cells.Item(132, 10).Value = 1012802   # J132
$ws.Cells.Item(132, 11).Value = 7394343   # K132
$ws.Cells.Item(132, 12).Value = 3038406   # L132
$ws.Cells.Item(132, 13).Value = -7391813   # M132
$ws.Cells.Item(132, 14).Value = -3043466   # N132
$ws.Cells.Item(136, 8).Value = 2456.2727   # H136
$ws.Cells.Item(136, 9).Value = 1651.9   # I136
$ws.Cells.Item(136, 11).Value = 4955.700000000001   # K136
$ws.Cells.Item(136, 13).Value = -2405.700000000001   # M136

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(86, 8).Value = 370208.3   # H86
$ws.Cells.Item(86, 9).Value = 1806   # I86
$ws.Cells.Item(86, 10).Value = 1168413.4   # J86
$ws.Cells.Item(86, 11).Value = 1806   # K86
$ws.Cells.Item(86, 12).Value = 1168413.4   # L86
$ws.Cells.Item(86, 13).Value = -683   # M86
$ws.Cells.Item(86, 14).Value = -1170659.4   # N86
$ws.Cells.Item(88, 8).Value = 8447.666999999999   # H88
$ws.Cells.Item(88, 10).Value = 8447.666999999999   # J88
$ws.Cells.Item(88, 12).Value = 8447.666999999999   # L88
$ws.Cells.Item(88, 14).Value = -9259.666999999999   # N88
$ws.Cells.Item(89, 8).Value = 370208.3   # H89
$ws.Cells.Item(89, 9).Value = 1806   # I89
$ws.Cells.Item(89, 10).Value = 1168413.4   # J89
$ws.Cells.Item(89, 11).Value = 9030   # K89
$ws.Cells.Item(89, 12).Value = 5842067   # L89
$ws.Cells.Item(89, 13).Value = -3414   # M89
$ws.Cells.Item(89, 14).Value = -5853299   # N89
$ws.Cells.Item(91, 8).Value = 8447.666999999999   # H91
$ws.Cells.Item(91, 10).Value = 8447.666999999999   # J91
$ws.Cells.Item(91, 12).Value = 8447.666999999999   # L91
$ws.Cells.Item(91, 14).Value = -11255.667   # N91
$ws.Cells.Item(107, 8).Value = 999.2857   # H107
$ws.Cells.Item(107, 9).Value = 998   # I107
$ws.Cells.Item(107, 10).Value = 999.8   # J107
$ws.Cells.Item(107, 11).Value = 998   # K107
$ws.Cells.Item(107, 12).Value = 999.8   # L107
$ws.Cells.Item(107, 13).Value = 922   # M107
$ws.Cells.Item(107, 14).Value = -4839.8   # N107
$ws.Cells.Item(134, 8).Value = 27387.2   # H134
$ws.Cells.Item(134, 9).Value = 1507.2413   # I134
$ws.Cells.Item(134, 10).Value = 95616.17999999999   # J134
$ws.Cells.Item(134, 11).Value = 4521.7239   # K134
$ws.Cells.Item(134, 12).Value = 286848.54   # L134
$ws.Cells.Item(134, 13).Value = -1986.7239   # M134
$ws.Cells.Item(134, 14).Value = -291918.54   # N134

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(99, 8).Value = 3630.75   # H99
$ws.Cells.Item(99, 9).Value = 3465.3635   # I99
$ws.Cells.Item(99, 10).Value = 3832.889   # J99
$ws.Cells.Item(99, 11).Value = 3465.3635   # K99
$ws.Cells.Item(99, 12).Value = 3832.889   # L99
$ws.Cells.Item(99, 13).Value = -1967.3635   # M99
$ws.Cells.Item(99, 14).Value = -6828.889   # N99
$ws.Cells.Item(126, 8).Value = 3630.75   # H126
$ws.Cells.Item(126, 9).Value = 3465.3635   # I126
$ws.Cells.Item(126, 10).Value = 3832.889   # J126
$ws.Cells.Item(126, 11).Value = 10396.0905   # K126
$ws.Cells.Item(126, 12).Value = 11498.667   # L126
$ws.Cells.Item(126, 13).Value = -7926.0905   # M126
$ws.Cells.Item(126, 14).Value = -16438.667   # N126
$ws.Cells.Item(134, 8).Value = 10205472   # H134
$ws.Cells.Item(134, 9).Value = 1154.1052   # I134
$ws.Cells.Item(134, 10).Value = 45456750   # J134
$ws.Cells.Item(134, 11).Value = 3462.3156   # K134
$ws.Cells.Item(134, 12).Value = 136370250   # L134
$ws.Cells.Item(134, 13).Value = -927.3155999999999   # M134
$ws.Cells.Item(134, 14).Value = -136375320   # N134

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(131, 8).Value = 28788580   # H131
$ws.Cells.Item(131, 9).Value = 363.07693   # I131
$ws.Cells.Item(131, 10).Value = 37699216   # J131
$ws.Cells.Item(131, 11).Value = 1089.23079   # K131
$ws.Cells.Item(131, 12).Value = 113097648   # L131
$ws.Cells.Item(131, 13).Value = 3950.76921   # M131
$ws.Cells.Item(131, 14).Value = -113107728   # N131

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(92, 8).Value = 8625.5   # H92
$ws.Cells.Item(92, 10).Value = 8625.5   # J92
$ws.Cells.Item(92, 12).Value = 8625.5   # L92
$ws.Cells.Item(92, 14).Value = -12369.5   # N92

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(64, 8).Value = 18571.428   # H64
$ws.Cells.Item(64, 9).Value = 10000   # I64
$ws.Cells.Item(64, 11).Value = 10000   # K64
$ws.Cells.Item(64, 13).Value = -9775   # M64
$ws.Cells.Item(67, 8).Value = 18571.428   # H67
$ws.Cells.Item(67, 9).Value = 10000   # I67
$ws.Cells.Item(67, 11).Value = 10000   # K67
$ws.Cells.Item(67, 13).Value = -9220   # M67
$ws.Cells.Item(136, 8).Value = 2526.5   # H136
$ws.Cells.Item(136, 9).Value = 1391.7142   # I136
$ws.Cells.Item(136, 11).Value = 4175.142599999999   # K136
$ws.Cells.Item(136, 13).Value = -1625.142599999999   # M136

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(69, 8).Value = 14000   # H69
$ws.Cells.Item(69, 10).Value = 14000   # J69
$ws.Cells.Item(69, 12).Value = 14000   # L69
$ws.Cells.Item(69, 14).Value = -15498   # N69
$ws.Cells.Item(72, 8).Value = 14000   # H72
$ws.Cells.Item(72, 10).Value = 14000   # J72
$ws.Cells.Item(72, 12).Value = 42000   # L72
$ws.Cells.Item(72, 14).Value = -49488   # N72
